$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R1")

$ws.Range("B4").Value = "R4"
$ws.Range("I4").Value = "SCECO+STB"
$ws.Range("J4").Value = "Good"
$ws.Range("L4").Value = "Latis"
